$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "initial budget" column (column B), shifting academic rank/in rank since left
$ws.Columns("B").Delete()

# Add headers for the new columns
$ws.Range("D1").Value = "DOB"
$ws.Range("E1").Value = "PATT yearly budget"
$ws.Range("F1").Value = "PA yearly budget"
$ws.Range("G1").Value = "PO yearly budget"

# Update row 2 (no longer has "initial budget"; academic rank now text "PA")
$ws.Range("B2").Value = "PA"

# Update row 3 values
$ws.Range("B3").Value = "PATT"

# Copy the date number format from C3 (existing date cell) onto D3 (new DOB cell)
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = 21916

$ws.Range("E3").Value = 12000
$ws.Range("F3").Value = 120000
$ws.Range("G3").Value = 240000

# Column widths for the new columns (D is a fixed manual width; E/F/G are
# sized to fit their header/content, same as the author's original sizing)
$ws.Columns("D").ColumnWidth = 18.17
$ws.Columns("E").ColumnWidth = 16.1
$ws.Columns("F").ColumnWidth = 13.99
$ws.Columns("G").ColumnWidth = 14.1
